# =====================================================================
#  "Json PATCH is now working" - applies the changes described by the
#  commit: rename Sheet1 -> ToDo, insert a new "Routing Rules" sheet
#  between ToDo and "Api Status Codes", fill in its content, mark the
#  JsonPatch checklist row as done ("Do!" -> "v", plus new "v" marks),
#  and refresh the selections/active-tab bookkeeping.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the first sheet "Sheet1" -> "ToDo"
# ---------------------------------------------------------------------
$todo = $wb.Worksheets.Item(1)
$todo.Name = "ToDo"

# ---------------------------------------------------------------------
# 2. Insert a brand-new "Routing Rules" sheet right after "ToDo" (i.e.
#    before "Api Status Codes")
# ---------------------------------------------------------------------
$statusCodes = $wb.Worksheets.Item("Api Status Codes")
$routing = $wb.Worksheets.Add($null, $todo)
$routing.Name = "Routing Rules"

# Column widths (B:F wide, G narrower) matching the checked-in sheet
$routing.Range("B:F").ColumnWidth = 33.42578125
$routing.Range("G:G").ColumnWidth = 16.42578125

# Header row (row 3)
$routing.Range("B3").Value = "GetById"
$routing.Range("C3").Value = "GetAll"
$routing.Range("D3").Value = "GetByRangeOfId"
$routing.Range("E3").Value = "Post"
$routing.Range("F3").Value = "Patch"
$routing.Range("G3").Value = "Delete"

# Example routes (row 4)
$routing.Range("B4").Value = "/pluralized/{id}"
$routing.Range("C4").Value = "/pluralized"
$routing.Range("D4").Value = "/range/pluralized"
$routing.Range("E4").Value = "/pluralized/new"
$routing.Range("F4").Value = "/pluralized/{id}"
$routing.Range("G4").Value = "/pluralized/{id}"

# Naming-convention note (row 11)
$routing.Range("B11").Value = 'all lowercase, separated by hyphens "-"'

# ---------------------------------------------------------------------
# 3. ToDo sheet: JsonPatch is done -> mark checklist cells with "v"
# ---------------------------------------------------------------------
# Row 10 (Controller): C..G = v
$todo.Range("C10:G10").Value = "v"
# Row 11 (Routing): C..E = v
$todo.Range("C11:E11").Value = "v"
# Row 13 (Status Codes): was "Do!" in C13, now "v" across C..E
$todo.Range("C13:E13").Value = "v"

# ---------------------------------------------------------------------
# 4. Selections / active sheet bookkeeping, per the saved workbook
# ---------------------------------------------------------------------
$todo.Activate()
$todo.Range("F13").Select()

$routing.Activate()
$routing.Range("F15").Select()

$statusCodes.Activate()
$statusCodes.Range("L16").Select()

# Leave "Routing Rules" as the active tab (matches activeTab="1")
$routing.Activate()
$routing.Range("F15").Select()
